# Fix shark double counts
# Updates the "Status by Landings" table values after correcting a double
# count in the shark landings figures. This adjusts the Sharks row (22),
# the Global totals row (23), and several rows whose aggregated totals
# are affected (6, 10, 11, 12, 15, 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (ISSCAAP 31)
$ws.Range("B6").Value = 1.161495750000001

# Row 10 (ISSCAAP 47)
$ws.Range("C10").Value = 1.39155272
$ws.Range("D10").Value = 7.735736368705187
$ws.Range("E10").Value = 39.52096861246186
$ws.Range("F10").Value = 52.74329501883294
$ws.Range("G10").Value = 47.25670498116705
$ws.Range("H10").Value = 52.74329501883294

# Row 11 (ISSCAAP 51)
$ws.Range("B11").Value = 5.480626450000003
$ws.Range("C11").Value = 5.225407434
$ws.Range("D11").Value = 30.90647844926512
$ws.Range("E11").Value = 44.26416256531046
$ws.Range("F11").Value = 24.82935898542442
$ws.Range("G11").Value = 75.17064101457558
$ws.Range("H11").Value = 24.82935898542442

# Row 12 (ISSCAAP 57)
$ws.Range("B12").Value = 5.88707478
$ws.Range("C12").Value = 5.79764829077466
$ws.Range("D12").Value = 17.57130444472443
$ws.Range("E12").Value = 65.67598477980499
$ws.Range("F12").Value = 16.75271077547058
$ws.Range("G12").Value = 83.24728922452942
$ws.Range("H12").Value = 16.75271077547058

# Row 15 (ISSCAAP 71)
$ws.Range("C15").Value = 11.463676695
$ws.Range("D15").Value = 38.90185099316687
$ws.Range("E15").Value = 21.06747471048451
$ws.Range("F15").Value = 40.03067429634861
$ws.Range("G15").Value = 59.96932570365138
$ws.Range("H15").Value = 40.03067429634861

# Row 16 (ISSCAAP 77)
$ws.Range("C16").Value = 1.708175780862069
$ws.Range("D16").Value = 51.87160332039483
$ws.Range("E16").Value = 33.63720569819994
$ws.Range("F16").Value = 14.49119098140523
$ws.Range("G16").Value = 85.50880901859477
$ws.Range("H16").Value = 14.49119098140523

# Row 22 (Sharks) - B22 picks up the same number format as C22 (s="9")
$ws.Range("B22").NumberFormat = $ws.Range("C22").NumberFormat
$ws.Range("B22").Value = 0.08399280000000002
$ws.Range("C22").Value = 0.05642497
$ws.Range("D22").Value = 48.68392486517937
$ws.Range("E22").Value = 37.54857113791996
$ws.Range("F22").Value = 13.76750399690066
$ws.Range("G22").Value = 86.23249600309933
$ws.Range("H22").Value = 13.76750399690066

# Row 23 (Global)
$ws.Range("B23").Value = 80.28049283
$ws.Range("C23").Value = 69.79157328945168
$ws.Range("D23").Value = 26.52088558783936
$ws.Range("E23").Value = 48.27323632967825
$ws.Range("F23").Value = 25.2058780824824
$ws.Range("G23").Value = 74.7941219175176
$ws.Range("H23").Value = 25.2058780824824
